$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two samples ("RM 232" and "SC 92") were dropped from the dataset, so remove those rows
# and let the remaining rows shift up to close the gap.
$ws.Rows(26).Delete()   # was "RM 232"
$ws.Rows(27).Delete()   # was "SC 92" (row 28 before the first delete)

# Refresh the randomly-missing value mask / recompute a couple of derived values
# for the remaining rows to match the newly regenerated dataset.
$ws.Range("D2").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("C6").Value = 15.1
$ws.Range("E6").Value = -5.7
$ws.Range("F7").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("F8").Value = 17.05
$ws.Range("F10").Value = 16.43
$ws.Range("E12").Value = ""
$ws.Range("E14").Value = -5.4
$ws.Range("F15").Value = 16.2
$ws.Range("C18").Value = 11.5
$ws.Range("F18").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E22").Value = -6.1
$ws.Range("C23").Value = 12.2
$ws.Range("F23").Value = ""
$ws.Range("F24").Value = ""
$ws.Range("C25").Value = ""
$ws.Range("E26").Value = ""
$ws.Range("B27").Value = -20.4
$ws.Range("E27").Value = ""
$ws.Range("B28").Value = ""
$ws.Range("E28").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("F29").Value = 18.06
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("F30").Value = 16.89
$ws.Range("E31").Value = -8.1
$ws.Range("B32").Value = ""
